$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "companiesPrimaryNames" column header in M1 ---------------
$ws.Cells.Item(1, 13).Value = "companiesPrimaryNames"

# New header cell uses an Arial 11 (black) font, matching the other header
# cells' style (non-wrapped alignment), but with its own font entry.
$ws.Cells.Item(1, 13).Font.Name = "Arial"
$ws.Cells.Item(1, 13).Font.Size = 11
$ws.Cells.Item(1, 13).Font.Color = 0

# --- Column M width ------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 20.14

# --- Extend the sheet's used range down to row 1000 (matches the rest of
# the pre-formatted rows already present on the sheet) --------------------
$ws.Cells.Item(1000, 13).NumberFormat = "General"

# --- View state: scroll so column F is left-most visible, and select H10 -
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 6
$aw.ScrollRow = 1
$ws.Range("H10").Select()
